# Add the new "23-jun" column (J) to the "Prix Spot" sheet, mirroring the
# style of the existing "22-jun" column (I) and filling in the new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the formatting used by the other date headers (e.g. I1)
# via Copy/PasteSpecial (xlPasteFormats), then set the text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "23-jun"

# New hourly values for 23-jun (rows 2-25 correspond to the 24 hourly slots).
$values = @(
    41.21,
    23.6,
    15.66,
    12.29,
    12.2,
    12.51,
    19.63,
    26,
    38.41,
    14.96,
    0.02,
    -0.01,
    -1.98,
    -8.890000000000001,
    -17.32,
    -9.529999999999999,
    -0.06,
    0,
    20.98,
    47.44,
    94.98,
    100.15,
    106.5,
    94.55
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
